# Insert a new weekly price record for "Haba" (Feria Lagunitas de Puerto Montt)
# above the existing row 79. This pushes all subsequent data rows (old 79-94)
# down by one (to 80-95), preserving each row's original "static" columns
# (Mercado/Region/Categoria/etc.) while the new row 79 carries the new
# date/volume/price/origin values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 79:94 down to 80:95, creating a blank row 79.
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with the new data record.
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 44782
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112026
$ws.Range("G79").Value = "Haba"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 80
$ws.Range("K79").Value = 17000
$ws.Range("L79").Value = 17000
$ws.Range("M79").Value = 17000
$ws.Range("N79").Value = "`$/saco 25 kilos"
$ws.Range("O79").Value = "Provincia de Limarí"
$ws.Range("P79").Value = 680
$ws.Range("Q79").Value = 25
$ws.Range("R79").Value = "Hortaliza"
